# Actualización automática de noticias - 2026-01-15
# Insert 3 new news rows at the top of the table (below the header row),
# pushing the existing rows down, then populate the new rows with the
# latest stories. The previously existing rows keep their content as-is
# (just shifted down), which Insert() handles for us automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 2 (old data rows shift
# down from 2-9 to 5-12).
$ws.Range("A2:F4").EntireRow.Insert()

# The inserted rows pick up the header row's bold/centered style by
# default; strip that back to the plain/default style used by the rest
# of the data rows.
$ws.Range("A2:F4").ClearFormats()

# Keep the "fecha" column as literal text (not auto-converted to a date
# serial) for the three new rows.
$ws.Range("A2:A4").NumberFormat = "@"

# Row 2 - new top story
$ws.Cells.Item(2, 1).Value = "2026-01-15"
$ws.Cells.Item(2, 2).Value = "La caída del “capo dei capi” de la Cosa Nostra: de una infancia dura al líder más temido de la mafia siciliana"
$ws.Cells.Item(2, 3).Value = "Infobae"
$ws.Cells.Item(2, 4).Value = "Sin identificar"
$ws.Cells.Item(2, 5).Value = "https://www.infobae.com/historias/2026/01/15/la-caida-del-capo-dei-capi-de-la-cosa-nostra-de-una-infancia-dura-al-lider-mas-temido-de-la-mafia-siciliana/"
$ws.Cells.Item(2, 6).Value = "PorGabriela Cicero"

# Row 3
$ws.Cells.Item(3, 1).Value = "2026-01-15"
$ws.Cells.Item(3, 2).Value = "El desplome de los colegios privados en Bogotá no se detiene: 35 cierran sus puertas para 2026"
$ws.Cells.Item(3, 3).Value = "Infobae"
$ws.Cells.Item(3, 4).Value = "Bogotá"
$ws.Cells.Item(3, 5).Value = "https://www.infobae.com/colombia/2026/01/14/el-desplome-de-los-colegios-privados-en-bogota-no-se-detiene-35-cierran-sus-puertas-para-2026/"
$ws.Cells.Item(3, 6).Value = "PorDaniella Mazo González"

# Row 4
$ws.Cells.Item(4, 1).Value = "2026-01-15"
$ws.Cells.Item(4, 2).Value = "Menor de edad colombiana tardó seis años para reencontrarse con su mamá en España, pero murió en grave accidente de tránsito"
$ws.Cells.Item(4, 3).Value = "Infobae"
$ws.Cells.Item(4, 4).Value = "Colombia"
$ws.Cells.Item(4, 5).Value = "https://www.infobae.com/colombia/2026/01/15/menor-de-edad-colombiana-tardo-seis-anos-para-reencontrarse-con-su-mama-en-espana-pero-murio-en-grave-accidente-de-transito/"
# (no resumen/F value for this story)
